$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns B through F (and their header/data cells) in a single pass,
# leaving only column A.
$ws.Columns("B:F").Delete()

# Column A becomes the lone, wider column (file-stored width 24).
# Excel's ColumnWidth property is expressed in "characters" and gets
# re-based by ~5/6 of a character when Excel converts it to the pixel grid
# before persisting the OOXML <col width> value, so back that offset out
# here to land on exactly width="24" in the saved file.
$ws.Columns.Item(1).ColumnWidth = 24 - (5/6)

# The header cell now describes the schedule-run button instead of the
# removed React-select inputs.
$ws.Range("A1").Value = "button_ScheduleTestRun"
